# Generate Report for Handback
# Updates status text, handback timestamps, clears the stale "handback not
# latest" error message (report regenerated), and widens/narrows a couple of
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Appears on every sheet (Overview E/F columns, and the Status column on
#    each per-locale sheet).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. Refreshed "Latest Handback DateTime" values for the newly generated
#    handback, plus clearing the stale "handback not latest" error text now
#    that the report is back in sync.
# ---------------------------------------------------------------------------
$wsZhCn.Range("K2").Value2 = "2016-11-08 23:59:22"
$wsZhCn.Range("P3").Value2 = ""

$wsDeDe.Range("K2").Value2 = "2016-11-08 23:59:40"
$wsDeDe.Range("P3").Value2 = ""

# ---------------------------------------------------------------------------
# 3. Column widths adjusted to fit the new (longer) status text in the
#    "Status" columns, and the now-empty "Error Detail" column shrinks back
#    down. ColumnWidth is quantized internally, so we pick the input that
#    lands the stored width as close as possible to the target.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
